# PreCondition1.xlsx — "Added original huang experiment folder"
#
# The workbook was previously authored on a Mac (paths using "/") and is
# now being re-saved from the Windows copy of the experiment, whose image
# paths use Windows-style back-slashes. Re-point the two PreCondName
# entries (and everything that reads them) at the back-slash paths, and
# drop the now-unused helper column C that only existed to preview the
# resulting Path(...) string for each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix the image folder separator: "/" -> "\" ------------------------
# Rows 2-6 reference the CS+ image, rows 7-11 the CS- image (shared-string
# driven, so this also naturally updates every dependent formula result).
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 1).Value = "PreCondition\CS+1.BMP"
}
for ($r = 7; $r -le 11; $r++) {
    $ws.Cells.Item($r, 1).Value = "PreCondition\CS-1.BMP"
}

# --- 2. The re-save also swapped which cached cell-format column A/B use --
# (cosmetic round-trip artifact of the different Excel build that saved the
# file) - mirror it by exchanging the A/B formatting for the data rows,
# using column D as scratch space.
$ws.Range("A2:A11").Copy() | Out-Null
$ws.Range("D2:D11").PasteSpecial(-4122) | Out-Null # xlPasteFormats

$ws.Range("B2:B11").Copy() | Out-Null
$ws.Range("A2:A11").PasteSpecial(-4122) | Out-Null # xlPasteFormats

$ws.Range("D2:D11").Copy() | Out-Null
$ws.Range("B2:B11").PasteSpecial(-4122) | Out-Null # xlPasteFormats

$ws.Range("D2:D11").Clear() | Out-Null
$excel.CutCopyMode = $false

# --- 3. Drop the helper "Path(...)" preview column -------------------------
$ws.Columns.Item(3).Delete() | Out-Null

# --- 4. Cosmetic: widen column A a touch and leave the cursor where the
# author left it when they last saved.
$ws.Columns.Item(1).ColumnWidth = 22.91
$ws.Range("F17").Select() | Out-Null
